$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.5185026666666667
$ws.Cells.Item(2, 8).Value = 1.555508
$ws.Cells.Item(2, 9).Value = 0.2217209591580833
$ws.Cells.Item(2, 10).Value = 0.2217209591580833
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.1528053333333333
$ws.Cells.Item(2, 14).Value = 0.458416
$ws.Cells.Item(2, 15).Value = 0.01103433215988526
$ws.Cells.Item(2, 16).Value = 0.01103433215988526
$ws.Cells.Item(2, 17).Value = 0.07922997281422221
$ws.Cells.Item(2, 18).Value = 0.713069755328
$ws.Cells.Item(2, 19).Value = 0.002446542710158645
$ws.Cells.Item(2, 20).Value = 0.002446542710158645
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.5185026666666667
$ws.Cells.Item(3, 8).Value = 1.555508
$ws.Cells.Item(3, 9).Value = 0.2217209591580833
$ws.Cells.Item(3, 10).Value = 0.2217209591580833
$ws.Cells.Item(3, 15).Value = 0.8539197603380489
$ws.Cells.Item(3, 16).Value = 0.8539197603380488
$ws.Cells.Item(3, 17).Value = 6.131412251941334
$ws.Cells.Item(3, 18).Value = 55.182710267472
$ws.Cells.Item(3, 19).Value = 0.1893319083061929
$ws.Cells.Item(3, 20).Value = 0.1893319083061928
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.5185026666666667
$ws.Cells.Item(4, 8).Value = 1.555508
$ws.Cells.Item(4, 9).Value = 0.2217209591580833
$ws.Cells.Item(4, 10).Value = 0.2217209591580833
$ws.Cells.Item(4, 13).Value = 1.712948333333333
$ws.Cells.Item(4, 14).Value = 5.138845
$ws.Cells.Item(4, 15).Value = 0.1236949029880405
$ws.Cells.Item(4, 16).Value = 0.1236949029880405
$ws.Cells.Item(4, 17).Value = 0.8881682786955556
$ws.Cells.Item(4, 18).Value = 7.993514508260001
$ws.Cells.Item(4, 19).Value = 0.02742575253347441
$ws.Cells.Item(4, 20).Value = 0.0274257525334744
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.5185026666666667
$ws.Cells.Item(5, 8).Value = 1.555508
$ws.Cells.Item(5, 9).Value = 0.2217209591580833
$ws.Cells.Item(5, 10).Value = 0.2217209591580833
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.1571906666666667
$ws.Cells.Item(5, 14).Value = 0.471572
$ws.Cells.Item(5, 15).Value = 0.01135100451402528
$ws.Cells.Item(5, 16).Value = 0.01135100451402528
$ws.Cells.Item(5, 17).Value = 0.08150377984177778
$ws.Cells.Item(5, 18).Value = 0.7335340185760001
$ws.Cells.Item(5, 19).Value = 0.002516755608257419
$ws.Cells.Item(5, 20).Value = 0.002516755608257418
$ws.Cells.Item(6, 9).Value = 0.4284141466794554
$ws.Cells.Item(6, 10).Value = 0.4284141466794553
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.1528053333333333
$ws.Cells.Item(6, 14).Value = 0.458416
$ws.Cells.Item(6, 15).Value = 0.01103433215988526
$ws.Cells.Item(6, 16).Value = 0.01103433215988526
$ws.Cells.Item(6, 17).Value = 0.1530899077991111
$ws.Cells.Item(6, 18).Value = 1.377809170192
$ws.Cells.Item(6, 19).Value = 0.004727263996454916
$ws.Cells.Item(6, 20).Value = 0.004727263996454915
$ws.Cells.Item(7, 9).Value = 0.4284141466794554
$ws.Cells.Item(7, 10).Value = 0.4284141466794553
$ws.Cells.Item(7, 15).Value = 0.8539197603380489
$ws.Cells.Item(7, 16).Value = 0.8539197603380488
$ws.Cells.Item(7, 19).Value = 0.3658313054579503
$ws.Cells.Item(7, 20).Value = 0.3658313054579502
$ws.Cells.Item(8, 9).Value = 0.4284141466794554
$ws.Cells.Item(8, 10).Value = 0.4284141466794553
$ws.Cells.Item(8, 13).Value = 1.712948333333333
$ws.Cells.Item(8, 14).Value = 5.138845
$ws.Cells.Item(8, 15).Value = 0.1236949029880405
$ws.Cells.Item(8, 16).Value = 0.1236949029880405
$ws.Cells.Item(8, 17).Value = 1.716138414112778
$ws.Cells.Item(8, 18).Value = 15.445245727015
$ws.Cells.Item(8, 19).Value = 0.05299264631221939
$ws.Cells.Item(8, 20).Value = 0.05299264631221937
$ws.Cells.Item(9, 9).Value = 0.4284141466794554
$ws.Cells.Item(9, 10).Value = 0.4284141466794553
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.1571906666666667
$ws.Cells.Item(9, 14).Value = 0.471572
$ws.Cells.Item(9, 15).Value = 0.01135100451402528
$ws.Cells.Item(9, 16).Value = 0.01135100451402528
$ws.Cells.Item(9, 17).Value = 0.1574834080848889
$ws.Cells.Item(9, 18).Value = 1.417350672764
$ws.Cells.Item(9, 19).Value = 0.004862930912830786
$ws.Cells.Item(9, 20).Value = 0.004862930912830785
$ws.Cells.Item(10, 7).Value = 0.6034673333333332
$ws.Cells.Item(10, 8).Value = 1.810402
$ws.Cells.Item(10, 9).Value = 0.2580533612824314
$ws.Cells.Item(10, 10).Value = 0.2580533612824314
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.1528053333333333
$ws.Cells.Item(10, 14).Value = 0.458416
$ws.Cells.Item(10, 15).Value = 0.01103433215988526
$ws.Cells.Item(10, 16).Value = 0.01103433215988526
$ws.Cells.Item(10, 17).Value = 0.09221302702577776
$ws.Cells.Item(10, 18).Value = 0.8299172432319999
$ws.Cells.Item(10, 19).Value = 0.002847446503365223
$ws.Cells.Item(10, 20).Value = 0.002847446503365223
$ws.Cells.Item(11, 7).Value = 0.6034673333333332
$ws.Cells.Item(11, 8).Value = 1.810402
$ws.Cells.Item(11, 9).Value = 0.2580533612824314
$ws.Cells.Item(11, 10).Value = 0.2580533612824314
$ws.Cells.Item(11, 15).Value = 0.8539197603380489
$ws.Cells.Item(11, 16).Value = 0.8539197603380488
$ws.Cells.Item(11, 17).Value = 7.136138807218666
$ws.Cells.Item(11, 18).Value = 64.225249264968
$ws.Cells.Item(11, 19).Value = 0.2203568644207218
$ws.Cells.Item(11, 20).Value = 0.2203568644207217
$ws.Cells.Item(12, 7).Value = 0.6034673333333332
$ws.Cells.Item(12, 8).Value = 1.810402
$ws.Cells.Item(12, 9).Value = 0.2580533612824314
$ws.Cells.Item(12, 10).Value = 0.2580533612824314
$ws.Cells.Item(12, 13).Value = 1.712948333333333
$ws.Cells.Item(12, 14).Value = 5.138845
$ws.Cells.Item(12, 15).Value = 0.1236949029880405
$ws.Cells.Item(12, 16).Value = 0.1236949029880405
$ws.Cells.Item(12, 17).Value = 1.033708362854444
$ws.Cells.Item(12, 18).Value = 9.303375265689999
$ws.Cells.Item(12, 19).Value = 0.03191988548956812
$ws.Cells.Item(12, 20).Value = 0.03191988548956812
$ws.Cells.Item(13, 7).Value = 0.6034673333333332
$ws.Cells.Item(13, 8).Value = 1.810402
$ws.Cells.Item(13, 9).Value = 0.2580533612824314
$ws.Cells.Item(13, 10).Value = 0.2580533612824314
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.1571906666666667
$ws.Cells.Item(13, 14).Value = 0.471572
$ws.Cells.Item(13, 15).Value = 0.01135100451402528
$ws.Cells.Item(13, 16).Value = 0.01135100451402528
$ws.Cells.Item(13, 17).Value = 0.09485943243822222
$ws.Cells.Item(13, 18).Value = 0.8537348919439999
$ws.Cells.Item(13, 19).Value = 0.002929164868776276
$ws.Cells.Item(13, 20).Value = 0.002929164868776275
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.2147046666666667
$ws.Cells.Item(14, 8).Value = 0.644114
$ws.Cells.Item(14, 9).Value = 0.09181153288002999
$ws.Cells.Item(14, 10).Value = 0.09181153288002998
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.1528053333333333
$ws.Cells.Item(14, 14).Value = 0.458416
$ws.Cells.Item(14, 15).Value = 0.01103433215988526
$ws.Cells.Item(14, 16).Value = 0.01103433215988526
$ws.Cells.Item(14, 17).Value = 0.03280801815822221
$ws.Cells.Item(14, 18).Value = 0.295272163424
$ws.Cells.Item(14, 19).Value = 0.001013078949906478
$ws.Cells.Item(14, 20).Value = 0.001013078949906478
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.2147046666666667
$ws.Cells.Item(15, 8).Value = 0.644114
$ws.Cells.Item(15, 9).Value = 0.09181153288002999
$ws.Cells.Item(15, 10).Value = 0.09181153288002998
$ws.Cells.Item(15, 15).Value = 0.8539197603380489
$ws.Cells.Item(15, 16).Value = 0.8539197603380488
$ws.Cells.Item(15, 17).Value = 2.538931635997333
$ws.Cells.Item(15, 18).Value = 22.850384723976
$ws.Cells.Item(15, 19).Value = 0.07839968215318412
$ws.Cells.Item(15, 20).Value = 0.07839968215318409
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.2147046666666667
$ws.Cells.Item(16, 8).Value = 0.644114
$ws.Cells.Item(16, 9).Value = 0.09181153288002999
$ws.Cells.Item(16, 10).Value = 0.09181153288002998
$ws.Cells.Item(16, 13).Value = 1.712948333333333
$ws.Cells.Item(16, 14).Value = 5.138845
$ws.Cells.Item(16, 15).Value = 0.1236949029880405
$ws.Cells.Item(16, 16).Value = 0.1236949029880405
$ws.Cells.Item(16, 17).Value = 0.3677780009255555
$ws.Cells.Item(16, 18).Value = 3.31000200833
$ws.Cells.Item(16, 19).Value = 0.0113566186527786
$ws.Cells.Item(16, 20).Value = 0.0113566186527786
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.2147046666666667
$ws.Cells.Item(17, 8).Value = 0.644114
$ws.Cells.Item(17, 9).Value = 0.09181153288002999
$ws.Cells.Item(17, 10).Value = 0.09181153288002998
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.1571906666666667
$ws.Cells.Item(17, 14).Value = 0.471572
$ws.Cells.Item(17, 15).Value = 0.01135100451402528
$ws.Cells.Item(17, 16).Value = 0.01135100451402528
$ws.Cells.Item(17, 17).Value = 0.03374956968977778
$ws.Cells.Item(17, 18).Value = 0.303746127208
$ws.Cells.Item(17, 19).Value = 0.001042153124160801
$ws.Cells.Item(17, 20).Value = 0.001042153124160801
